$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (Changed) date column C, rows 2-20, from serial 45181 to 45182
for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 45182
}
